$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 84; this pushes the existing rows 84-89 down to 85-90
# and shifts any formatting (like the date style on column D) along with them.
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new weekly record.
$ws.Cells.Item(84, 1).Value = 11
$ws.Cells.Item(84, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(84, 3).Value = 'Bíobío'
$ws.Cells.Item(84, 4).Value = 45147
$ws.Cells.Item(84, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(84, 5).Value = 8
$ws.Cells.Item(84, 6).Value = 100112031
$ws.Cells.Item(84, 7).Value = 'Poroto verde'
$ws.Cells.Item(84, 8).Value = 'Magnum'
$ws.Cells.Item(84, 9).Value = 'Primera'
$ws.Cells.Item(84, 10).Value = 50
$ws.Cells.Item(84, 11).Value = 28000
$ws.Cells.Item(84, 12).Value = 28000
$ws.Cells.Item(84, 13).Value = 28000
$ws.Cells.Item(84, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(84, 15).Value = 'Perú'
$ws.Cells.Item(84, 16).Value = 1120
$ws.Cells.Item(84, 17).Value = 25
$ws.Cells.Item(84, 18).Value = 'Hortaliza'
